# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 633
    $ws.Range("F3").Value = 580
    $ws.Range("F6").Value = 129
    $ws.Range("F8").Value = 59
    $ws.Range("F10").Value = 5088
    $ws.Range("F11").Value = 4736
}
